$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the "center" alignment style (matches cellXfs index 3) to D4 and D7:D21,
# creating empty but styled cells in column D for the Status checklist.
$ws.Range("D4").HorizontalAlignment = -4108
$ws.Range("D7:D21").HorizontalAlignment = -4108

# D6 gets the same style plus an "x" mark (same as D5).
$ws.Range("D6").HorizontalAlignment = -4108
$ws.Range("D6").Value = "x"

# Update the active selection from D5 to E5.
$ws.Range("E5").Select()
